$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue "D2" "21.664.84"
Set-TextValue "E2" "  -2.05%  "
Set-TextValue "D3" "1.537.26"
Set-TextValue "E3" "  -1.32%  "
Set-TextValue "E4" "  +0.21%  "
Set-TextValue "E5" "  +0.20%  "
Set-TextValue "D6" "290.14"
Set-TextValue "E6" "  -0.74%  "
Set-TextValue "D7" "0.3931"
Set-TextValue "E7" "  -0.68%  "
Set-TextValue "D8" "0.3197"
Set-TextValue "D9" "43.61"
Set-TextValue "E9" "  +0.49%  "
Set-TextValue "D10" "0.07199"
Set-TextValue "E10" "  -1.95%  "
Set-TextValue "D11" "1.078"
Set-TextValue "E11" "  -1.41%  "
Set-TextValue "E12" "  +0.13%  "
Set-TextValue "D13" "5.780"
Set-TextValue "E13" "  +1.97%  "
Set-TextValue "D14" "18.36"
Set-TextValue "E14" "  -4.12%  "
Set-TextValue "D15" "6.622"
Set-TextValue "E15" "  -1.02%  "
Set-TextValue "D16" "1.538.68"
Set-TextValue "E16" "  -1.82%  "
Set-TextValue "D17" "0.00001094"
Set-TextValue "E17" "  -4.67%  "
Set-TextValue "D18" "0.06619"
Set-TextValue "E18" "  +0.18%  "
Set-TextValue "D19" "84.22"
Set-TextValue "E19" "  -0.13%  "
Set-TextValue "D20" "1.001"
Set-TextValue "E20" "  +0.24%  "
Set-TextValue "D21" "6.148"
Set-TextValue "E21" "  -2.94%  "
Set-TextValue "D22" "15.51"
Set-TextValue "E22" "  -2.42%  "
Set-TextValue "E23" "  -4.63%  "
Set-TextValue "D24" "2.386"
Set-TextValue "E24" "  +1.89%  "
Set-TextValue "D25" "21.671.99"
Set-TextValue "E25" "  -2.06%  "
Set-TextValue "D26" "2.396"
Set-TextValue "E26" "  -2.53%  "
Set-TextValue "D27" "152.62"
Set-TextValue "E27" "  +2.98%  "
Set-TextValue "E28" "  -1.10%  "
Set-TextValue "D29" "4.900"
Set-TextValue "E29" "  +0.67%  "
Set-TextValue "D30" "1.710.88"
Set-TextValue "E30" "  -1.66%  "
Set-TextValue "D31" "117.71"
Set-TextValue "E31" "  -1.36%  "
Set-TextValue "D32" "6.047"
Set-TextValue "E32" "  +5.44%  "
Set-TextValue "D33" "0.9596"
Set-TextValue "E33" "  -8.65%  "
Set-TextValue "D34" "0.08078"
Set-TextValue "E34" "  -3.95%  "
Set-TextValue "D35" "5.202"
Set-TextValue "E35" "  +0.60%  "
Set-TextValue "D36" "8.513"
Set-TextValue "E36" "  -6.82%  "
Set-TextValue "D37" "1.486"
Set-TextValue "E37" "  -8.63%  "
Set-TextValue "D38" "0.02219"
Set-TextValue "E38" "  -2.86%  "
Set-TextValue "B39" "Aptos"
Set-TextValue "C39" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D39" "11.39"
Set-TextValue "E39" "  +5.31%  "
Set-TextValue "B40" "Hedera"
Set-TextValue "C40" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D40" "0.05906"
Set-TextValue "E40" "  -4.61%  "
Set-TextValue "D41" "0.2044"
Set-TextValue "E41" "  -1.69%  "
Set-TextValue "D42" "1.183"
Set-TextValue "E42" "  -2.98%  "
Set-TextValue "D43" "1.001"
Set-TextValue "E43" "  +0.26%  "
Set-TextValue "D44" "0.5812"
Set-TextValue "E44" "  -1.16%  "
Set-TextValue "D45" "13.20"
Set-TextValue "E45" "  +0.13%  "
Set-TextValue "D46" "3.730"
Set-TextValue "E46" "  -1.12%  "
Set-TextValue "D47" "0.5586"
Set-TextValue "E47" "  -0.83%  "
Set-TextValue "D48" "1.894"
Set-TextValue "E48" "  -1.03%  "
Set-TextValue "D49" "1.161"
Set-TextValue "E49" "  +1.28%  "
Set-TextValue "D50" "116.02"
Set-TextValue "E50" "  -2.64%  "
Set-TextValue "E51" "  -2.35%  "
